$d = $word.ActiveDocument

$newText = "Isso é uma ação especial, devido ao formato da aula com muitas informações sobre software"

# The document ends with a paragraph whose trailing (collapsed) range carries
# the auto-managed "_GoBack" bookmark. We are appending a new paragraph after
# it, and the bookmark needs to end up collapsed at the end of that new
# paragraph (mirroring what Word does when you type at the end of a doc).

# 1) Drop the existing _GoBack bookmark - it will be re-created below at the
#    correct (new) location.
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# 2) Append a new paragraph with the requested sentence at the very end of
#    the document.
$r = $d.Content
$r.Collapse(0)
$r.InsertAfter("`r" + $newText)

# 3) Re-create "_GoBack" collapsed right after the new sentence (i.e. right
#    before the final paragraph mark). Placing a bookmark directly on a
#    zero-width range located exactly at a paragraph-mark boundary is
#    unreliable, so two throw-away characters are inserted there first, the
#    bookmark is attached around the first (now an interior, unambiguous
#    position), and both placeholder characters are removed again - leaving
#    the bookmark correctly collapsed in place without touching the document
#    text.
$insertPos = $d.Content.End - 1

$mr = $d.Range($insertPos, $insertPos)
$mr.InsertAfter("ZW")

$bmRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$zRange = $d.Bookmarks("_GoBack").Range
$zRange.Delete()

$wRange = $d.Range($insertPos, $insertPos + 1)
$wRange.Delete()
